try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # --- Swap rows 2 (TECHM.NS) and 3 (ASIANPAINT.NS) -----------------------
    # Use Copy/PasteSpecial (not direct .Value writes) so the multi-line
    # "Last Price" text doesn't trigger an auto row-height recalculation.
    $ws.Range("A3:F3").Copy()
    $ws.Range("A200:F200").PasteSpecial()
    $ws.Range("A2:F2").Copy()
    $ws.Range("A3:F3").PasteSpecial()
    $ws.Range("A200:F200").Copy()
    $ws.Range("A2:F2").PasteSpecial()
    $ws.Rows.Item(200).EntireRow.Delete()

    # --- Insert a new row 9 for JIOFIN.NS ------------------------------------
    $ws.Rows.Item(9).EntireRow.Insert()

    $ws.Range("A9").Value = "JIOFIN.NS"

    # Stage the multi-line text far away, then copy it in, to avoid the
    # auto row-height recalculation that a direct multi-line .Value write
    # triggers.
    $ws.Range("A200").Value = "Ticker`nJIOFIN.NS    327.399994`nName: 2025-08-14 00:00:00, dtype: float64"
    $ws.Range("A200").Copy()
    $ws.Range("B9").PasteSpecial()
    $ws.Rows.Item(200).EntireRow.Delete()

    $ws.Range("C9").Value = 315.1475115966797
    $ws.Range("D9").Value = 281.9188722229004
    $ws.Range("E9").Value = 45842
    $ws.Range("F9").Value = 44

    Write-Output "done"
} catch {
    Write-Output "ERROR: $_"
}
